$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Length (in)"
$ws.Range("K1").Value = "Width (in)"
$ws.Range("L1").Value = "Height (in)"

$ws.Range("J2").Value = 4.0599999999999996
$ws.Range("K2").Value = 1.36
$ws.Range("L2").Value = 1.91

$ws.Range("J3").Value = 5.44
$ws.Range("K3").Value = 1.69
$ws.Range("L3").Value = 1.64

$ws.Range("J4").Value = 5.46
$ws.Range("K4").Value = 1.69
$ws.Range("L4").Value = 1.99

$ws.Range("J5").Value = 5.43
$ws.Range("K5").Value = 1.65
$ws.Range("L5").Value = 2.04

$ws.Range("J1:L1").Font.Bold = $true
